# Apply the weekly update described in the commit "Fruta / hortaliza, semanal".
# Two new observation rows are inserted into the date-ordered data table,
# pushing the subsequent rows down (one insertion point after row 54,
# another insertion point once the table has shifted, at row 139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new record (new row 55) ---------------------------------
$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value = 5
$ws.Range("B55").Value = 'Macroferia Regional de Talca'
$ws.Range("C55").Value = 'Maule'
$ws.Range("D55").Value = 44965
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 100112030
$ws.Range("G55").Value = 'Poroto granado'
$ws.Range("H55").Value = 'Sin especificar'
$ws.Range("I55").Value = 'Primera'
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = 25000
$ws.Range("N55").Value = '$/saco 25 kilos'
$ws.Range("O55").Value = 'Región del Maule'
$ws.Range("P55").Value = 1000
$ws.Range("Q55").Value = 25
$ws.Range("R55").Value = 'Hortaliza'

# --- Insert second new record (new row 139, after the first shift) --------
$ws.Rows.Item(139).Insert()

$ws.Range("A139").Value = 5
$ws.Range("B139").Value = 'Macroferia Regional de Talca'
$ws.Range("C139").Value = 'Maule'
$ws.Range("D139").Value = 44964
$ws.Range("E139").Value = 7
$ws.Range("F139").Value = 100112030
$ws.Range("G139").Value = 'Poroto granado'
$ws.Range("H139").Value = 'Sin especificar'
$ws.Range("I139").Value = 'Primera'
$ws.Range("J139").Value = 500
$ws.Range("K139").Value = 28000
$ws.Range("L139").Value = 28000
$ws.Range("M139").Value = 28000
$ws.Range("N139").Value = '$/saco 25 kilos'
$ws.Range("O139").Value = 'Región del Maule'
$ws.Range("P139").Value = 1120
$ws.Range("Q139").Value = 25
$ws.Range("R139").Value = 'Hortaliza'
